# Updated cryptos list on Sat Aug  5 13:17:07 UTC 2023 with GitHub Actions
# Refreshes Price (column D) and Volume(1h) (column E) text values for each
# coin row, and reorders the Aave / Mantle / BabyDogeCoin rows (46-48).
#
# Price values are entered as literal text (they use '.' as a thousands
# separator, e.g. "29.000.99", and some carry significant trailing zeros,
# e.g. "1.000"), matching how the source data is stored in the workbook.
# For values that look like plain numbers (e.g. "241.41") a leading
# apostrophe forces Excel to keep them as text instead of auto-converting
# to a numeric value; the cell style is then reset to "Normal" so no stray
# text-format style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '28.997.40'
$ws.Range("E2").Value = '  -0.70%  '
$ws.Range("D3").Value = '1.830.22'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("D4").Value = "'0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'241.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").Value = "'0.6259"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.48%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = "'45.30"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.22%  '
$ws.Range("D9").Value = "'0.07579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.89%  '
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("E11").Value = '  -2.29%  '
$ws.Range("D12").Value = "'0.07748"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.16%  '
$ws.Range("D13").Value = '1.827.09'
$ws.Range("E13").Value = '  -1.22%  '
$ws.Range("D14").Value = "'4.952"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").Value = "'0.6621"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.43%  '
$ws.Range("D16").Value = "'82.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").Value = "'0.000009466"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +9.32%  '
$ws.Range("D18").Value = "'5.975"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.81%  '
$ws.Range("D19").Value = '28.979.79'
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("D20").Value = "'223.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.75%  '
$ws.Range("D21").Value = "'12.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.65%  '
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = "'7.196"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").Value = "'158.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.62%  '
$ws.Range("D26").Value = "'8.409"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("D27").Value = "'0.1364"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.26%  '
$ws.Range("D28").Value = "'17.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.34%  '
$ws.Range("D29").Value = "'1.489"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").Value = "'4.049"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("D31").Value = "'4.014"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").Value = "'1.194"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("D33").Value = "'0.05188"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.67%  '
$ws.Range("E34").Value = '  -1.99%  '
$ws.Range("D35").Value = "'0.7361"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.75%  '
$ws.Range("E36").Value = '  -0.88%  '
$ws.Range("D37").Value = "'2.698"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.80%  '
$ws.Range("D38").Value = '1.255.57'
$ws.Range("E38").Value = '  -4.61%  '
$ws.Range("D39").Value = "'2.763"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").Value = "'6.223"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.88%  '
$ws.Range("D42").Value = "'0.8883"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.41%  '
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").Value = "'101.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").Value = '1.974.65'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'64.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").Value = "'0.5110"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").Value = "'0.3972"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.06%  '
$ws.Range("D50").Value = "'8.851"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.39%  '
$ws.Range("D51").Value = "'0.05749"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.61%  '
